# Rencita_Day_10_Assingment_2 - Question_2 completed - Annotations pending
$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# --- Constants sheet: four new Name/Value rows (18-21) --------------------
$wsConstants.Range("A18").Value = "OutputWorkbookPath"
$wsConstants.Range("B18").Value = "Data/Output/Assignment1.xlsx"

$wsConstants.Range("A19").Value = "OutputSheetName"
$wsConstants.Range("B19").Value = "NSE"

$wsConstants.Range("A20").Value = "BrowserPath"
$wsConstants.Hyperlinks.Add($wsConstants.Range("B20"), "https://www.nseindia.com/")

$wsConstants.Range("A21").Value = "NSEBrowserPath"
$wsConstants.Hyperlinks.Add($wsConstants.Range("B21"), "https://www.nseindia.com/get-quotes/equity?symbol=")

# --- Settings sheet: new Name/Value pair in row 6 --------------------------
$wsSettings.Range("A6").Value = "MailCredentials"
$wsSettings.Range("B6").Value = "MailCredentails"

# --- Assets sheet: new Asset row (row 2) ------------------------------------
$wsAssets.Range("A2").Value = "ToMailAddress"
$wsAssets.Range("B2").Value = "ToMailAddress"
$wsAssets.Range("C2").Value = "RPATraining"

# --- Selections left active on each sheet after editing --------------------
$wsConstants.Activate()
$wsConstants.Range("A3").Select()
$wsConstants.Range("B26").Select()

$wsAssets.Activate()
$wsAssets.Range("A2").Select()

$wsSettings.Activate()
$wsSettings.Range("A6").Select()
